$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header column C: "Nhận tin nhắn" (bold, matching existing header style)
$ws.Range("C1").Value = "Nhận tin nhắn"
$ws.Range("C1").Font.Bold = $true

# Fill data rows for new column
$ws.Range("C2").Value = "Có"
$ws.Range("C3").Value = "Không"

# Update the active selection to match the authored state
$ws.Range("E6").Select()
